# CRUD in MVC Core
# Add Assignment_15 row to the Assignments tracking sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newUrl = "https://github.com/Vasanth30e/Assignments_Phase2/tree/master/Assignment_15/CustomErrorHandling"

# Fill in the new row's values first.
$ws.Range("A16").Value = "Assignment_15"
$ws.Range("B16").Value = $newUrl
$ws.Range("C16").Value = 45169

# Turn B16 into a hyperlink (adds the external relationship + hyperlink entry).
$ws.Hyperlinks.Add($ws.Range("B16"), $newUrl)

# Row 14 (the row just above the former last row) carries the exact cell
# formatting we want for the new row - copy it across so the new row matches
# the rest of the table (borders, fonts, alignment) instead of inheriting the
# "closing" border style that sits on the former last row (row 15).
$ws.Range("A14:C14").Copy()
$ws.Range("A16:C16").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Row height for the new row
$ws.Rows.Item(16).RowHeight = 36

# Update selection to the newly added cell
$ws.Range("B16").Select()

# Update window view position
$excel.ActiveWindow.Left = 7980
$excel.ActiveWindow.Top = 3855
